$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 4040
$ws.Range("B7").Value = "dddd"
$ws.Range("C7").Value = "zzzz"
$ws.Range("D7").Value = "2nd line added"

$ws.Range("A8").Select()
